$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.731.04"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.72%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.627.19"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.12%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.60%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "211.20"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.521"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.63%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.21"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E9").Value = "  -2.78%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0612"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.34%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0879"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.83%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.863.34"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.83%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.631.92"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.563"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.17"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.56%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "27.800.58"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.51%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "229.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.58%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0721"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.49%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.50"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.04%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.55%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.36"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.85%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.35"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.98%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.78%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "154.44"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.94"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "15.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -1.03%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -1.11%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0481"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.92%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.41"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.50%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.403.03"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.65%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.07"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +7.59%  "
$ws.Range("E37").Value = "  +1.10%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0169"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.558"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.866"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("E41").Value = "  -0.14%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "66.49"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "5.51"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.772.10"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "87.82"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0997"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0506"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0₇0978"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -7.61%  "
